# edit.ps1 - apply xlsx diff: add evaluator_partial_correctness column,
# refresh o_10 data, add o_20 and o_20_jumbled sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "o_10": add column E header + values, refresh row 2 ----
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("E1").Value = 'evaluator_partial_correctness'
$ws1.Range("D1").Copy() | Out-Null
$ws1.Range("E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$s1_a2 = @'
 Given is the adjacency matrix for a weighted directed graph containing 16 nodes labelled A to P. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the least cost path from node A to node P? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P
 A 0 2 0 0 3 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 2 0 0 3 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 1 0 0 4 0 0 0 0 0 0 0
 F 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 0 0 3 0 0 4 0 0 0 0 0
 H 0 0 0 0 0 0 3 0 0 0 0 3 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0
 J 0 0 0 0 0 4 0 0 4 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 3 0 0 0 0 2 0
 L 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 2
 M 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
Solution: A -> B -> C -> D -> H -> L -> P
 Given these examples, answer the following quesiton.
what is the least cost path from node A to node P? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P
 A 0 5 0 0 2 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 2 0 0 3 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 2 0 0 2 0 0 0 0 0 0 0
 F 0 0 0 0 3 0 0 0 0 2 0 0 0 0 0 0
 G 0 0 4 0 0 5 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 0 0 0 2 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 1 0 0 3 0 0 0
 J 0 0 0 0 0 0 0 0 0 0 4 0 0 4 0 0
 K 0 0 0 0 0 0 3 0 0 0 0 0 0 0 5 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 3
 M 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 4 0 3 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
    
'@
$ws1.Range("A2").Value = $s1_a2
$ws1.Range("B2").Value = 'A -> E -> I -> J -> N -> O -> P'
$ws1.Range("C2").Value = 'The least cost path from node A to node P is A -> E -> F -> J -> K -> L -> P.'
$ws1.Range("D2").Value = 'invalid input'
$ws1.Range("E2").Value = '2/7'

# ---- New sheet "o_20", inserted right after "o_10" ----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "o_20"

$ws1.Range("A1:E1").Copy() | Out-Null
$ws2.Range("A1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws2.Range("A1").Value = 'prompt'
$ws2.Range("B1").Value = 'solution'
$ws2.Range("C1").Value = 'llm_response'
$ws2.Range("D1").Value = 'evaluator_response'
$ws2.Range("E1").Value = 'evaluator_partial_correctness'

$s2_a2 = @'
 Given is the adjacency matrix for a weighted directed graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the least cost path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 5 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 2 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 5 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 3 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 4 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 2 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 3 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 5 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 3 0 0 0 4 0 0 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 5 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 2 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 3 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 2 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 2 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 5 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 2 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 1 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0
Solution: A -> B -> G -> H -> I -> J -> O -> T -> Y
 Given these examples, answer the following quesiton.
what is the least cost path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 2 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 3 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 5 0 5 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 2 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 5 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 5 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 5 0 2 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 5 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 5 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 3 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 5 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 5 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 4 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
    
'@
$ws2.Range("A2").Value = $s2_a2
$ws2.Range("B2").Value = 'A -> B -> G -> H -> I -> N -> O -> T -> Y'
$ws2.Range("C2").Value = 'Solution: A -> B -> G -> H -> I -> J -> T -> Y'
$ws2.Range("D2").Value = 'invalid input'
$ws2.Range("E2").Value = '5/9'

# ---- New sheet "o_20_jumbled", inserted right after "o_20" ----
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "o_20_jumbled"

$ws1.Range("A1:E1").Copy() | Out-Null
$ws3.Range("A1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws3.Range("A1").Value = 'prompt'
$ws3.Range("B1").Value = 'solution'
$ws3.Range("C1").Value = 'llm_response'
$ws3.Range("D1").Value = 'evaluator_response'
$ws3.Range("E1").Value = 'evaluator_partial_correctness'

$s3_a2 = @'
 Given is the adjacency matrix for a weighted directed graph containing 24 nodes labelled A to X. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the least cost path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 5 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 2 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 5 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 3 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 4 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 2 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 3 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 5 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 3 0 0 0 4 0 0 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 5 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 2 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 3 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 2 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 2 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 5 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 2 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 1 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0
Solution: A -> B -> G -> H -> I -> J -> O -> T -> Y
 Given these examples, answer the following quesiton.
what is the least cost path from node A to node X? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X
 A 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 3 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 5 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 0 0 3 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 4 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 2 0 0 0 0 0 3 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 5 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 4 0 0 0 4 0 0 0 0 0 3 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 2 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 2 0 0 0 2 0 1 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 2 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 5 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 4 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
    
'@
$ws3.Range("A2").Value = $s3_a2
$ws3.Range("B2").Value = 'A -> F -> G -> H -> M -> R -> S -> W -> X'
$ws3.Range("C2").Value = 'The least cost path from node A to node X is: A -> F -> G -> H -> I -> N -> X.'
$ws3.Range("D2").Value = 'invalid input'
$ws3.Range("E2").Value = '4/9'

$ws1.Select()
$ws1.Range("A1").Select()
